# Apply the "Updated cryptos list" snapshot refresh (GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Excel auto-converts numeric-looking strings (e.g. "0.0950") to real
    # numbers on assignment, which would drop significant trailing zeros and
    # change the stored cell type away from text. Temporarily mark the cell
    # as Text so the literal string is kept, then restore its original
    # (Normal) style so no stray formatting is left behind.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "54.000.00"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "2.252.74"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "492.31"
$ws.Range("E5").Value = "  +1.30%  "
Set-TextValue "D6" "127.30"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +0.67%  "
Set-TextValue "D9" "0.0950"
$ws.Range("E9").Value = "  +2.99%  "
Set-TextValue "D10" "0.151"
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("E11").Value = "  +3.18%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "2.655.08"
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("D15").Value = "53.919.60"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "2.251.57"
$ws.Range("E17").Value = "  +1.56%  "
Set-TextValue "D18" "9.99"
$ws.Range("E18").Value = "  +4.64%  "
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D20" "6.42"
$ws.Range("E20").Value = "  +4.63%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "299.55"
$ws.Range("E21").Value = "  +1.84%  "
Set-TextValue "D22" "0.999"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  -1.67%  "
Set-TextValue "D24" "61.94"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("E25").Value = "  +1.89%  "
Set-TextValue "D26" "0.368"
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("D27").Value = "2.355.30"
$ws.Range("E27").Value = "  +2.49%  "
Set-TextValue "D28" "0.148"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("E29").Value = "  +0.59%  "
Set-TextValue "D30" "166.65"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D32" "5.84"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0₃0674"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +0.29%  "
Set-TextValue "D36" "1.06"
$ws.Range("E36").Value = "  -0.15%  "
Set-TextValue "D37" "17.59"
$ws.Range("E37").Value = "  +1.88%  "
Set-TextValue "D38" "0.884"
$ws.Range("E38").Value = "  +7.20%  "
$ws.Range("E39").Value = "  +2.68%  "
Set-TextValue "D40" "3.64"
$ws.Range("E40").Value = "  +3.15%  "
Set-TextValue "D41" "35.70"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("E42").Value = "  +1.35%  "
Set-TextValue "D43" "1.38"
$ws.Range("E43").Value = "  +1.84%  "
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("E45").Value = "  +2.92%  "
Set-TextValue "D46" "124.73"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("E47").Value = "  +0.54%  "
Set-TextValue "D48" "0.539"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  +2.29%  "
Set-TextValue "D50" "234.57"
$ws.Range("E50").Value = "  +1.34%  "
Set-TextValue "D51" "0.0201"
$ws.Range("E51").Value = "  +0.76%  "
